$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9520253539085388
$ws.Range("B1").Value = 2.08130145072937
$ws.Range("C1").Value = 8.01541805267334
$ws.Range("D1").Value = 2.42518949508667
$ws.Range("E1").Value = 0.776460587978363
